$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# "Fix 0 hospitalisation error for India": several rows in the "type"
# column (I) were wrongly tagged as date/timestamp/integer instead of
# character, which broke downstream parsing. Correct them.
$rows = @(2, 3, 4, 12, 39, 56, 58, 60)
foreach ($r in $rows) {
    $ws.Cells.Item($r, 9).Value = "character"
}

# Update the view state to reflect where the fix was made.
$ws.Range("I2:I69").Select()
$excel.ActiveWindow.ScrollRow = 24
